$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A41").Value = "Don't forget to delete m, p, e from updated SBTs."
$ws.Range("A42").Value = "Why did it work?"

[void]$ws.Range("A43").Select()

